# Pseudonymizer FAKE_DATA_2019 update:
# - new function creates a pseudonymous id for people sharing address+last name
# - a handful of existing cell values change (civil-status codes, house number suffix)
# - a new synthetic row (11) is appended to the data table
# - active selection moves to C8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing-cell value updates -----------------------------------------

# Row 2: civil status code "Civil stat#B" -> "Civil stat#L"
$ws.Range("F2").Value = "Civil stat#L"

# Row 9: birthdate becomes an (invalid, hence textual) date "31.02.1965"
$ws.Range("D9").Value = "31.02.1965"

# --- new row 11 -------------------------------------------------------
# Copy formatting from row 9 (same per-column styles: A -> integer, D -> date)
# so the new row reuses the existing style indices instead of creating new ones.
$ws.Range("A9:T9").Copy()
$ws.Range("A11:T11").PasteSpecial(-4122)

$ws.Range("A11").Value = 7560000000009
$ws.Range("B11").Value = "Apmann"
$ws.Range("C11").Value = "Ali"
$ws.Range("D11").Value = 30246
$ws.Range("E11").Value = "Sex#W"
$ws.Range("F11").Value = "Civil stat#L"
$ws.Range("G11").Value = "Antragssteller"
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = "CHResidenc#L"
$ws.Range("J11").Value = "Ackerstrasse"
$ws.Range("K11").Value = 11
$ws.Range("L11").Value = 804500
$ws.Range("M11").Value = "Zürich"
$ws.Range("N11").Value = 77777
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = -30

# Row 9: civil status code "Civil stat#B" -> "Civil stat#M"
$ws.Range("F9").Value = "Civil stat#M"

# Row 4: house-number suffix "3c" -> "3c!"
$ws.Range("K4").Value = "3c!"

# Row 10: birthdate serial 34940 -> 26320 (still a real date)
$ws.Range("D10").Value = 26320
# Row 10: civil status code "Civil stat#B" -> "Civil stat#M"
$ws.Range("F10").Value = "Civil stat#M"

# --- selection -------------------------------------------------------
$ws.Range("C8").Select()
